{"js": "// el-287: inseridas as testemunhas\n// The witness-name template placeholders referenced the wrong field\n// (\"name.first\"); point them at \"name.text\" instead, for both witnesses.\n\nconst body = context.document.body;\n\nconst replacements = [\n  { oldText: \"witnesses[0].name.first\", newText: \"witnesses[0].name.text\" },\n  { oldText: \"witnesses[1].name.first\", newText: \"witnesses[1].name.text\" },\n];\n\nfor (const { oldText, newText } of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# el-287: inseridas as testemunhas\n# The witness-name template placeholders referenced the wrong field\n# (\"name.first\"); point them at \"name.text\" instead, for both witnesses.\n\n$replacements = @(\n    @{ Old = \"witnesses[0].name.first\"; New = \"witnesses[0].name.text\" },\n    @{ Old = \"witnesses[1].name.first\"; New = \"witnesses[1].name.text\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Forward = $true\n    $find.Wrap = 1            # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2) | Out-Null\n}\n"}
